$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sheet is protected; unprotect before editing
$ws.Unprotect()

# Update the confidential disclaimer date text in A10
$ws.Range("A10").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-04-09 for illustrative purposes only and are subject to change."

# Update Weight (D) and Percent Change (E) values for rows 2-6
$ws.Range("D2").Value = 0.489259866700984
$ws.Range("E2").Value = 0

$ws.Range("D3").Value = 0.3324772788614663
$ws.Range("E3").Value = 0.004015670910871849

$ws.Range("D4").Value = 0.09401945441713461
$ws.Range("E4").Value = -0.0007459903021259651

$ws.Range("D5").Value = 0.05464470608274609
$ws.Range("E5").Value = -0.0006876002750401122

$ws.Range("D6").Value = 0.02959869393766897
$ws.Range("E6").Value = -0.01650269755633149

# Update Percent Change total in E7 (Total row)
$ws.Range("E7").Value = 0.0007389497269956191

# Re-protect the sheet to restore original protection state
$ws.Protect()
